# Update "想去人数" (interested-attendee counts) on the "展览" and "全部类型"
# sheets to reflect the latest generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5451   # 南宁·AB动漫游戏嘉年华
$ws1.Range("F4").Value = 11830  # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）
$ws1.Range("F6").Value = 599    # 南宁·蔚蓝档案only
$ws1.Range("F7").Value = 174    # 南宁·国乙only
$ws1.Range("F8").Value = 297    # 南宁·熊喵M动漫嘉年华【免费】
$ws1.Range("F9").Value = 1063   # 南宁·第二届北极光动漫展

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5451   # 南宁·AB动漫游戏嘉年华
$ws4.Range("F7").Value = 11830  # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）
$ws4.Range("F9").Value = 599    # 南宁·蔚蓝档案only
$ws4.Range("F10").Value = 174   # 南宁·国乙only
$ws4.Range("F13").Value = 297   # 南宁·熊喵M动漫嘉年华【免费】
$ws4.Range("F14").Value = 1063  # 南宁·第二届北极光动漫展
